# Auto-generated edit script applying the crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.710.83'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '''3.343.07'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.70%  '
$ws.Range("D5").Value = '''576.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '''178.08'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("D7").Value = '''0.618'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.55%  '
$ws.Range("D9").Value = '''3.340.59'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.70%  '
$ws.Range("D10").Value = '''0.130'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = '''6.89'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '''0.414'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '''3.920.07'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.80%  '
$ws.Range("D14").Value = '''0.135'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = '''28.79'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.80%  '
$ws.Range("D16").Value = '''65.708.47'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").Value = '''0.0000170'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '''3.348.69'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.66%  '
$ws.Range("E19").Value = '  -3.03%  '
$ws.Range("D20").Value = '''13.45'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("D21").Value = '''364.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''71.53'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '''0.999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000123'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '''0.522'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("D27").Value = '''9.63'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("D31").Value = '''5.66'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '''22.89'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.01%  '
$ws.Range("D34").Value = '''6.88'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D35").Value = '''1.23'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("E36").Value = '  -2.16%  '
$ws.Range("D37").Value = '''160.21'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("E38").Value = '  -4.46%  '
$ws.Range("D39").Value = '''27.36'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.34%  '
$ws.Range("D40").Value = '''1.76'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = '''2.56'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '''2.701.24'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("D43").Value = '''4.30'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").Value = '''6.24'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.79%  '
$ws.Range("D45").Value = '''0.0669'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("D46").Value = '''336.56'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +8.69%  '
$ws.Range("D47").Value = '''39.66'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("D48").Value = '''24.38'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("E49").Value = '  -3.15%  '
$ws.Range("D50").Value = '''0.104'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.26%  '
$ws.Range("D51").Value = '''0.970'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.34%  '
